# Revision: added pyrolysis parameter row ("revision, added pyrolysis and
# additional figures").
#
# The target state inserts one new parameter row, "chemical_recycling_pyrolysis"
# (value TRUE), immediately after the existing "chemical_recycling_gasification"
# row (row 9). Every row that used to follow it (old rows 10-24, i.e.
# fossil_routes ... fossil_lock_in, together with their value/explanation
# figures) simply shifts down by one row to rows 11-25, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10. This pushes the existing rows 10:24
# down to 11:25, carrying their values/formatting along untouched.
$ws.Rows.Item(10).Insert()

# Populate the newly-inserted row for the new pyrolysis parameter.
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true
